$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "meta"
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("meta")

$meta.Range("A2").Value = 45988
$meta.Range("B2").Value = 0.45833333333333331
$meta.Range("C2").Value = ""

# Row 2 no longer needs the taller (wrapped) row height.
$meta.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------
# Sheet "events"
# ---------------------------------------------------------------------
$events = $wb.Worksheets.Item("events")

# --- Row 2 (id 1) ---
$events.Range("C2").Value = "San Martín"
$events.Range("D2").Value = "Chipurana"
$events.Range("E2").Value = 45987
$events.Range("F2").Value = "Incendio urbano"
$events.Range("G2").Value = "05 viviendas inhabitables"
$events.Range("H2").Value = "Sin novedad"
$events.Range("I2").Value = "Activa"
$events.Range("K2").Value = -6.4340320000000002
$events.Range("L2").Value = -75.666836000000004

# --- Row 3 (id 2) ---
$events.Range("B3").Value = "San Martín"
$events.Range("C3").Value = "Huallaga"
$events.Range("D3").Value = "Alto Saposoa"
$events.Range("E3").Value = 45988
$events.Range("F3").Value = "Lluvias intensas"
$events.Range("G3").Value = "Viviendas y vías afectadas"
$events.Range("H3").Value = "En monitoreo"
$events.Range("K3").Value = -6.6018619999999997
$events.Range("L3").Value = -76.963115000000002

# --- Row 4 (id 3) ---
# This row previously carried explicit cell styles on E:I - clear them so the
# cells fall back to the default (unstyled) formatting, then set the values.
$events.Range("E4:I4").ClearFormats()
$events.Range("B4").Value = "San Martín"
$events.Range("C4").Value = "Bellavista"
$events.Range("D4").Value = "Alto Biavo"
$events.Range("E4").Value = 45988
$events.Range("F4").Value = "Lluvias intensas"
$events.Range("G4").Value = "Viviendas y vías afectadas"
$events.Range("H4").Value = "En monitoreo"
$events.Range("I4").Value = "Activa"
$events.Range("K4").Value = -7.8221109999999996
$events.Range("L4").Value = -76.274123000000003

# --- Row 5 (id 4) ---
# Entirely removed from the data set; clear contents but keep the
# pre-existing per-cell styles on E5:I5.
$events.Range("A5:L5").ClearContents()

# ---------------------------------------------------------------------
# Sheet selections / active sheet
# ---------------------------------------------------------------------
[void]$events.Range("B5").Select()
$meta.Activate()
[void]$meta.Range("C2").Select()
